$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.676.63"
$ws.Range("E2").Value = "  -2.18%  "

$ws.Range("D3").Value = "2.962.94"
$ws.Range("E3").Value = "  -3.32%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'495.52"
$ws.Range("E5").Value = "  -6.03%  "

$ws.Range("D6").Value = "'134.70"
$ws.Range("E6").Value = "  -5.96%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.425"
$ws.Range("E8").Value = "  -5.34%  "

$ws.Range("D9").Value = "'7.16"
$ws.Range("E9").Value = "  -6.67%  "

$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -6.47%  "

$ws.Range("D11").Value = "'0.351"
$ws.Range("E11").Value = "  -5.45%  "

$ws.Range("D12").Value = "3.493.37"
$ws.Range("E12").Value = "  -2.81%  "

$ws.Range("E13").Value = "  -3.07%  "

$ws.Range("D14").Value = "'25.63"
$ws.Range("E14").Value = "  -6.70%  "

$ws.Range("D15").Value = "'0.0000156"
$ws.Range("E15").Value = "  -8.27%  "

$ws.Range("D16").Value = "56.842.45"
$ws.Range("E16").Value = "  -1.93%  "

$ws.Range("D17").Value = "2.981.39"
$ws.Range("E17").Value = "  -2.82%  "

$ws.Range("D18").Value = "'5.99"
$ws.Range("E18").Value = "  -3.84%  "

$ws.Range("D19").Value = "'12.45"
$ws.Range("E19").Value = "  -5.83%  "

$ws.Range("D20").Value = "'7.73"
$ws.Range("E20").Value = "  -5.64%  "

$ws.Range("D21").Value = "'315.79"
$ws.Range("E21").Value = "  -7.52%  "

$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").Value = "'5.70"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").Value = "'0.485"
$ws.Range("E24").Value = "  -3.61%  "

$ws.Range("D25").Value = "'62.80"
$ws.Range("E25").Value = "  -3.32%  "

$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -5.76%  "

$ws.Range("D28").Value = "0.0₃0864"
$ws.Range("E28").Value = "  -11.81%  "

$ws.Range("D29").Value = "'6.49"
$ws.Range("E29").Value = "  -6.81%  "

$ws.Range("D30").Value = "'6.95"
$ws.Range("E30").Value = "  -6.31%  "

$ws.Range("D31").Value = "'1.75"
$ws.Range("E31").Value = "  -5.94%  "

$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  -8.79%  "

$ws.Range("D33").Value = "'19.87"
$ws.Range("E33").Value = "  -5.80%  "

$ws.Range("D34").Value = "'153.90"
$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").Value = "'4.47"
$ws.Range("E35").Value = "  -6.57%  "

$ws.Range("D36").Value = "'5.67"
$ws.Range("E36").Value = "  -5.76%  "

$ws.Range("D37").Value = "'1.21"
$ws.Range("E37").Value = "  -8.80%  "

$ws.Range("D38").Value = "'23.74"
$ws.Range("E38").Value = "  -10.45%  "

$ws.Range("D39").Value = "'0.0652"
$ws.Range("E39").Value = "  -7.73%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'37.61"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.001.59"
$ws.Range("E41").Value = "  -3.20%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").Value = "'0.640"
$ws.Range("E43").Value = "  -4.22%  "

$ws.Range("D44").Value = "'3.65"
$ws.Range("E44").Value = "  -6.56%  "

$ws.Range("D45").Value = "2.153.76"
$ws.Range("E45").Value = "  -7.86%  "

$ws.Range("D46").Value = "'1.35"
$ws.Range("E46").Value = "  -8.88%  "

$ws.Range("D47").Value = "'5.85"
$ws.Range("E47").Value = "  -3.20%  "

$ws.Range("D48").Value = "'0.919"
$ws.Range("E48").Value = "  -10.91%  "

$ws.Range("D49").Value = "'0.0230"
$ws.Range("E49").Value = "  -6.05%  "

$ws.Range("D50").Value = "'18.90"
$ws.Range("E50").Value = "  -6.53%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0857"
$ws.Range("E51").Value = "  -4.83%  "
